# Applies the "Updated symbol list" commit: refreshes the Price column (D)
# for most coins, and re-ranks three rows (41-43) so that KickToken now
# appears above BKEXToken and CEJI, with updated prices/volume labels.
#
# All Price cells in this sheet are stored as text (not numbers), so each
# one is forced to Text format ("@") before assignment. This preserves the
# exact textual representation (e.g. trailing zeros like "0.05712" or
# "0.00000000750") instead of letting Excel auto-convert the numeric-looking
# string into a floating point value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Row 2: BNB ---
Set-TextValue "D2" '249.06'

# --- Row 3: OKB ---
Set-TextValue "D3" '22.55'

# --- Row 4: HuobiToken ---
Set-TextValue "D4" '5.396'

# --- Row 5: Cronos ---
Set-TextValue "D5" '0.05712'

# --- Row 6: GateToken ---
Set-TextValue "D6" '3.410'

# --- Row 7: KuCoinToken ---
Set-TextValue "D7" '6.335'

# --- Row 8: MXToken ---
Set-TextValue "D8" '0.8142'

# --- Row 9: FTXToken ---
Set-TextValue "D9" '0.9283'

# --- Row 10: WazirX ---
Set-TextValue "D10" '0.1422'

# --- Row 11: MandalaExchangeToken ---
Set-TextValue "D11" '0.07524'

# --- Row 12: LiechtensteinCryptoassetsExchange ---
Set-TextValue "D12" '0.03119'

# --- Row 13: BitrueCoin ---
Set-TextValue "D13" '0.03047'

# --- Row 14: BitMartToken ---
Set-TextValue "D14" '0.09346'

# --- Row 16: BitForexToken ---
Set-TextValue "D16" '0.001598'

# --- Row 17: CoinExToken ---
Set-TextValue "D17" '0.04770'

# --- Row 19: One ---
Set-TextValue "D19" '0.0005794'

# --- Row 20: TigerCash ---
Set-TextValue "D20" '0.006453'

# --- Row 21: HotbitToken ---
Set-TextValue "D21" '0.005007'

# --- Row 22: BitKan ---
Set-TextValue "D22" '0.001025'

# --- Row 25: BTSEToken ---
Set-TextValue "D25" '2.167'

# --- Row 27: ProBitToken ---
Set-TextValue "D27" '0.1308'

# --- Row 40: IDEX ---
Set-TextValue "D40" '0.03991'

# --- Row 41: now KickToken (was BKEXToken) ---
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D41" '0.006863'
$ws.Range("E41").Value = '40KickTokenKICK'

# --- Row 42: now BKEXToken (was CEJI) ---
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1066'
$ws.Range("E42").Value = '41BKEXTokenBKK'

# --- Row 43: now CEJI (was KickToken) ---
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002711'
$ws.Range("E43").Value = '42CEJICEJI'

# --- Row 44: LocalTraders ---
Set-TextValue "D44" '0.007535'

# --- Row 45: CoinLion ---
Set-TextValue "D45" '0.00005894'

# --- Row 46: Kangarootoken ---
Set-TextValue "D46" '0.00000000750'

# --- Row 47: CoinbaseStockToken ---
Set-TextValue "D47" '0.5003'

# --- Row 49: CryptobidCoin ---
Set-TextValue "D49" '0.00002101'
